$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NegativeLogin")
$ws.Activate()

# Update the "username and password do not match" expected-error message text
# used by the negative-login test rows (row 2 and row 4).
$newText = "Epic sadface: Username and password do not match any user in this service"
$ws.Range("C2").Value = $newText
$ws.Range("C4").Value = $newText

# Move the active selection to C4, matching where the cursor was left after edits.
$ws.Range("C4").Select()
